$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Linear Regression: an algorithm which attempts to fit an equation of the form Y = wX + c to data so that its error amount to the true value is minimised as far as possible",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Linear Regression: an algorithm which attempts to fit an equation of the form Y = wX + c to the data so that its error to the true value is minimised as far as possible",
    2
)
